# The deck's overall look was switched from the "Integral" theme to the
# built-in "Office Theme" (this is what PowerPoint's Design gallery does
# when you pick the plain "Office Theme" design for the whole deck).
# The visible, rendering-relevant effect of that switch is the slide
# master's 12-slot theme colour scheme (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink) changing from the Integral palette to the Office
# palette, so recreate that via the Design's ColorScheme.

function ConvertTo-ComColor([string]$hex) {
    # PowerPoint/VBA .RGB values are packed little-endian (0x00BBGGRR),
    # i.e. r + g*256 + b*65536, given a "RRGGBB" hex string.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.ColorScheme

# Office Theme palette, in the same dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink
# order that ColorScheme.Colors(1..12) addresses.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = ConvertTo-ComColor $officeThemeColors[$i]
}
